# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx data refresh described by the diff:
#  - Updates Price (column D) and Volume(1h) (column E) text values for rows 2-47
#  - Inserts a new 'BabyDogeCoin' entry at row 48, shifting Cronos and Mantle
#    down one row each and dropping USDD from the list (sheet stays 51 rows)
#  - Updates row 51 (EnergySwap) price/volume
#
# A handful of Price cells (D5, D11, D16, D29, D49) are numeric-looking text
# that carries a significant trailing zero (e.g. "212.70", "0.0510"). Excel's
# automatic type inference would otherwise coerce those into Number cells and
# silently drop the trailing zero, so those specific cells are pre-formatted
# as Text before the value is written, exactly as typing them in the UI would
# require to keep the trailing zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.485.56'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '1.626.07'
$ws.Range('E3').Value = '  -0.76%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.70'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('E6').Value = '  +1.06%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -0.36%  '
$ws.Range('E9').Value = '  -2.11%  '
$ws.Range('D10').Value = '18.73'
$ws.Range('E10').Value = '  -1.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0840'
$ws.Range('E11').Value = '  +0.37%  '
$ws.Range('D12').Value = '1.851.87'
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('D13').Value = '1.616.15'
$ws.Range('E13').Value = '  -1.36%  '
$ws.Range('E14').Value = '  +1.10%  '
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.90'
$ws.Range('E16').Value = '  +2.76%  '
$ws.Range('D17').Value = '26.517.26'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('E18').Value = '  -0.53%  '
$ws.Range('D19').Value = '213.96'
$ws.Range('E19').Value = '  +2.35%  '
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('E21').Value = '  -0.86%  '
$ws.Range('E22').Value = '  +1.35%  '
$ws.Range('D23').Value = '9.27'
$ws.Range('E23').Value = '  -1.49%  '
$ws.Range('D24').Value = '2.04'
$ws.Range('E24').Value = '  +5.53%  '
$ws.Range('D25').Value = '148.56'
$ws.Range('E25').Value = '  +1.38%  '
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('E27').Value = '  -0.98%  '
$ws.Range('D28').Value = '6.83'
$ws.Range('E28').Value = '  +0.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.50'
$ws.Range('E29').Value = '  +0.56%  '
$ws.Range('E31').Value = '  -1.09%  '
$ws.Range('D32').Value = '3.32'
$ws.Range('E32').Value = '  +2.50%  '
$ws.Range('E33').Value = '  -1.10%  '
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('D35').Value = '1.225.27'
$ws.Range('E35').Value = '  +4.79%  '
$ws.Range('E36').Value = '  -1.18%  '
$ws.Range('D37').Value = '0.0173'
$ws.Range('E37').Value = '  +3.16%  '
$ws.Range('E38').Value = '  +0.18%  '
$ws.Range('D39').Value = '0.794'
$ws.Range('E39').Value = '  -2.12%  '
$ws.Range('D40').Value = '0.505'
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  -1.81%  '
$ws.Range('D42').Value = '0.793'
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('E43').Value = '  -0.76%  '
$ws.Range('D44').Value = '1.761.23'
$ws.Range('E44').Value = '  -0.85%  '
$ws.Range('D45').Value = '92.83'
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('D47').Value = '54.75'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0101'
$ws.Range('E48').Value = '  -3.64%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0510'
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '0.406'
$ws.Range('E50').Value = '  -0.73%  '
$ws.Range('D51').Value = '7.47'
$ws.Range('E51').Value = '  -0.99%  '
